# Applies the "Aplicando Funcoes no Crud" edit:
#   - Colors the body/list text (but not the bold section headers, the
#     intro paragraph, or the blank/hr paragraphs) with accent6 green
#     (RGB 4EA72E / theme color accent6), including the paragraph mark.
#
# Strategy: Word's Font.Color setter correctly produces
#   <w:color w:val="4EA72E"/>
# on every run (and the paragraph-mark rPr) in the target range, but it
# cannot also stamp the w:themeColor="accent6" attribute (the COM shim's
# ObjectThemeColor / RGB setters each clobber the other instead of
# merging). So for each target paragraph we: set Font.Color via the
# normal object model (gets the structure + w:val exactly right), then
# read back that single paragraph's WordOpenXML, patch the resulting
# w:val="4EA72E" colors to also carry w:themeColor="accent6", and feed
# it back with Range.InsertXML scoped to that same paragraph so only its
# content changes.

$d = $word.ActiveDocument

$targetParagraphs = @(6,7,10,11,12,13,16,17,18,19,20,23,24,25,26,27)

foreach ($idx in $targetParagraphs) {
    $p = $d.Paragraphs($idx)
    $r = $p.Range

    # 1) Apply the RGB color through the normal object model so every
    #    run - and the paragraph mark - gets a <w:color w:val="4EA72E"/>
    #    in exactly the right place.
    $r.Font.Color = 3057486

    # 2) Patch in the accompanying theme attribute by round-tripping this
    #    paragraph's own OOXML (InsertXML replaces only this range).
    $xml = $r.WordOpenXML
    $xml = $xml -replace '<w:color w:val="4EA72E"/>', '<w:color w:val="4EA72E" w:themeColor="accent6"/>'
    $r.InsertXML($xml)
}

Write-Output "done"
